# Handback status report refresh: a new round of handoff/handback xliff
# generation completed for the "49fb647f-...md" source file, so its
# timestamps move forward on the Overview sheet and on each per-locale
# (zh-cn / de-de) detail sheet. The "d2cf8647-...md" file's row (row 3)
# was not touched in this run, so it is left as-is.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 = Latest HO Xliff Generate Date for 49fb647f-...md
$wsOverview.Range("G2").Value = "2016-08-29 06:49:10"

# zh-cn sheet, row 2 (49fb647f-...md): Correspond Handoff / Handback Datetime
$wsZhCn.Range("H2").Value = "2016-08-29 06:49:00"
$wsZhCn.Range("K2").Value = "2016-08-29 06:49:27"

# de-de sheet, row 2 (49fb647f-...md): Correspond Handoff / Handback Datetime
$wsDeDe.Range("H2").Value = "2016-08-29 06:49:10"
$wsDeDe.Range("K2").Value = "2016-08-29 06:49:34"
